$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 244 - this shifts the existing rows 244:395 down to 245:396,
# carrying all their original values/formatting with them.
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with its data.
$ws.Range("A244").Value = 5
$ws.Range("B244").Value = "Macroferia Regional de Talca"
$ws.Range("C244").Value = "Maule"
$ws.Range("D244").Value = 44824
$ws.Range("E244").Value = 7
$ws.Range("F244").Value = 100114014
$ws.Range("G244").Value = "Betarraga"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Segunda"
$ws.Range("J244").Value = 3000
$ws.Range("K244").Value = 800
$ws.Range("L244").Value = 800
$ws.Range("M244").Value = 800
$ws.Range("N244").Value = "`$/paquete 5 unidades"
$ws.Range("O244").Value = "Región del Maule"
$ws.Range("P244").Value = 160
$ws.Range("Q244").Value = 5
$ws.Range("R244").Value = "Hortaliza"
